# Update "想去人数" (want-to-go count, column F) and "最低票价" (min price, column G)
# figures on the "展览" (sheet1) and "全部类型" (sheet4) worksheets to reflect the
# refreshed bilibili scrape output (gh-pages build 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1434
$ws1.Range("G3").Value = 65
$ws1.Range("G5").Value = 25
$ws1.Range("F7").Value = 11985
$ws1.Range("F8").Value = 4453
$ws1.Range("F10").Value = 54
$ws1.Range("F14").Value = 1114
$ws1.Range("F15").Value = 174
$ws1.Range("F16").Value = 57
$ws1.Range("F17").Value = 5190
$ws1.Range("F21").Value = 11400
$ws1.Range("F22").Value = 11419

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1434
$ws4.Range("G3").Value = 65
$ws4.Range("G5").Value = 25
$ws4.Range("F7").Value = 11985
$ws4.Range("F8").Value = 4453
$ws4.Range("F10").Value = 54
$ws4.Range("F15").Value = 1114
$ws4.Range("F16").Value = 174
$ws4.Range("F17").Value = 57
$ws4.Range("F18").Value = 5190
$ws4.Range("F22").Value = 11400
$ws4.Range("F23").Value = 11419
